# Update Lgi1-Adam22 LR-pair sheet with new TPM-derived values.
# Sending cluster changes from "Neutrophils" to "MuSCs" for every data row,
# and the per-row NATMI-derived numeric columns (F,G,H,M,N,O,P,Q,R,S,T) are
# recomputed accordingly. Target cluster (column D) text is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: Sending cluster -> "MuSCs" for all data rows (2-7)
$ws.Range("A2:A7").Value = "MuSCs"

# Row 2 (Target cluster: ECs)
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.035285
$ws.Range("H2").Value = 0.07056999999999999
$ws.Range("M2").Value = 5.447678
$ws.Range("N2").Value = 10.895356
$ws.Range("O2").Value = 0.4286498436662743
$ws.Range("P2").Value = 0.4047900009176674
$ws.Range("Q2").Value = 0.19222131823
$ws.Range("R2").Value = 0.7688852729199999
$ws.Range("S2").Value = 0.4286498436662743
$ws.Range("T2").Value = 0.4047900009176674

# Row 3 (Target cluster: FAPs)
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.035285
$ws.Range("H3").Value = 0.07056999999999999
$ws.Range("M3").Value = 0.7327576666666668
$ws.Range("N3").Value = 2.198273
$ws.Range("O3").Value = 0.05765694287766837
$ws.Range("P3").Value = 0.08167139556406268
$ws.Range("Q3").Value = 0.02585535426833334
$ws.Range("R3").Value = 0.15513212561
$ws.Range("S3").Value = 0.05765694287766837
$ws.Range("T3").Value = 0.08167139556406268

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("F4").Value = 0.5
$ws.Range("G4").Value = 0.035285
$ws.Range("H4").Value = 0.07056999999999999
$ws.Range("M4").Value = 0.302684
$ws.Range("N4").Value = 0.9080520000000001
$ws.Range("O4").Value = 0.02381665165971311
$ws.Range("P4").Value = 0.03373642586009028
$ws.Range("Q4").Value = 0.01068020494
$ws.Range("R4").Value = 0.06408122964
$ws.Range("S4").Value = 0.02381665165971311
$ws.Range("T4").Value = 0.03373642586009028

# Row 5 (Target cluster: MuSCs)
$ws.Range("F5").Value = 0.5
$ws.Range("G5").Value = 0.035285
$ws.Range("H5").Value = 0.07056999999999999
$ws.Range("M5").Value = 5.763022
$ws.Range("N5").Value = 11.526044
$ws.Range("O5").Value = 0.4534626457997884
$ws.Range("P5").Value = 0.4282216534583244
$ws.Range("Q5").Value = 0.20334823127
$ws.Range("R5").Value = 0.8133929250799999
$ws.Range("S5").Value = 0.4534626457997884
$ws.Range("T5").Value = 0.4282216534583244

# Row 6 (Target cluster: Neutrophils)
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.035285
$ws.Range("H6").Value = 0.07056999999999999
$ws.Range("M6").Value = 0.2632226666666667
$ws.Range("N6").Value = 0.789668
$ws.Range("O6").Value = 0.02071164171525676
$ws.Range("P6").Value = 0.02933816118029118
$ws.Range("Q6").Value = 0.009287811793333333
$ws.Range("R6").Value = 0.05572687076
$ws.Range("S6").Value = 0.02071164171525676
$ws.Range("T6").Value = 0.02933816118029118

# Row 7 (Target cluster: Resolving-Mac)
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.035285
$ws.Range("H7").Value = 0.07056999999999999
$ws.Range("M7").Value = 0.199559
$ws.Range("N7").Value = 0.598677
$ws.Range("O7").Value = 0.01570227428129894
$ws.Range("P7").Value = 0.02224236301956415
$ws.Range("Q7").Value = 0.007041439315
$ws.Range("R7").Value = 0.04224863589
$ws.Range("S7").Value = 0.01570227428129894
$ws.Range("T7").Value = 0.02224236301956415
